# Fix NLA location / crosswalk data entry issue on the "data" sheet:
# The surface-sample columns (AB:AI = sample.depth.s ... turb.s) and the
# depth-sample columns (AX:BE = sample.depth.d ... turb.d) had been entered
# swapped for rows 3-17.  Swap the two blocks back into the correct columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

for ($r = 3; $r -le 17; $r++) {
    $surfaceRange = $ws.Range("AB$r`:AI$r")
    $depthRange   = $ws.Range("AX$r`:BE$r")

    $surfaceValues = $surfaceRange.Value()
    $depthValues   = $depthRange.Value()

    $surfaceRange.Value = $depthValues
    $depthRange.Value   = $surfaceValues
}

# Correct the mis-keyed longitude value in D4 (data sheet)
$ws.Range("D4").Value = 41.676720000000003

# Leave the view focused on the corrected area, matching where editing left off
$ws.Activate() | Out-Null
$ws.Range("AC11").Select() | Out-Null

